$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap B4 (Out of PO) and B7 (NSY) values
$ws.Range("B4").Value = 399
$ws.Range("B7").Value = 401

# Add new transaction row describing the trade, matching style of previous note rows
$ws.Range("A15").Value = "08.03.2025 - Out of PO Bam Adebayo karşılığında NSY'ye 2 Dolara vermiştir. (399-401)"
$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial(-4122)

# Update selection to reflect new active cell after edit
$ws.Range("A16").Select()
